# Applies the "cryptos list" refresh described in the commit:
#   "Updated cryptos list on Sun May 21 16:41:47 UTC 2023 with GitHub Actions"
#
# For each affected row in the crypto table (columns: B=Coin, C=Link, D=Price, E=Volume(1h))
# this updates the Price / Volume(1h) text, and for rows 47-48 also swaps the
# Coin name + Link (PaxDollar now ranks above EnergySwap).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: Row number, then any of B/C/D/E values that changed ($null = unchanged)
$updates = @(
    [pscustomobject]@{ Row = 2; B = $null; C = $null; D = "27.106.15"; E = "  -0.32%  " }
    [pscustomobject]@{ Row = 3; B = $null; C = $null; D = "1.824.48"; E = "  -0.52%  " }
    [pscustomobject]@{ Row = 4; B = $null; C = $null; D = "1.012"; E = "  +0.02%  " }
    [pscustomobject]@{ Row = 5; B = $null; C = $null; D = "312.24"; E = "  -0.45%  " }
    [pscustomobject]@{ Row = 6; B = $null; C = $null; D = $null; E = "  -0.15%  " }
    [pscustomobject]@{ Row = 7; B = $null; C = $null; D = "0.4630"; E = "  -1.69%  " }
    [pscustomobject]@{ Row = 8; B = $null; C = $null; D = "0.3628"; E = "  -1.54%  " }
    [pscustomobject]@{ Row = 9; B = $null; C = $null; D = "0.07288"; E = "  -1.71%  " }
    [pscustomobject]@{ Row = 10; B = $null; C = $null; D = "0.8695"; E = "  -1.44%  " }
    [pscustomobject]@{ Row = 11; B = $null; C = $null; D = "20.11"; E = "  -1.52%  " }
    [pscustomobject]@{ Row = 12; B = $null; C = $null; D = "1.877.61"; E = "  +3.66%  " }
    [pscustomobject]@{ Row = 13; B = $null; C = $null; D = "0.07634"; E = "  +4.20%  " }
    [pscustomobject]@{ Row = 14; B = $null; C = $null; D = "5.339"; E = "  -2.59%  " }
    [pscustomobject]@{ Row = 15; B = $null; C = $null; D = "92.29"; E = "  -0.51%  " }
    [pscustomobject]@{ Row = 16; B = $null; C = $null; D = "6.470"; E = "  -1.37%  " }
    [pscustomobject]@{ Row = 17; B = $null; C = $null; D = $null; E = "  -0.25%  " }
    [pscustomobject]@{ Row = 18; B = $null; C = $null; D = "0.000008603"; E = "  -2.13%  " }
    [pscustomobject]@{ Row = 19; B = $null; C = $null; D = "1.009"; E = "  -0.05%  " }
    [pscustomobject]@{ Row = 20; B = $null; C = $null; D = "27.449.76"; E = "  +0.91%  " }
    [pscustomobject]@{ Row = 21; B = $null; C = $null; D = $null; E = "  -2.03%  " }
    [pscustomobject]@{ Row = 22; B = $null; C = $null; D = "5.213"; E = "  -1.65%  " }
    [pscustomobject]@{ Row = 23; B = $null; C = $null; D = "10.57"; E = "  -1.11%  " }
    [pscustomobject]@{ Row = 24; B = $null; C = $null; D = "2.098.78"; E = "  +2.62%  " }
    [pscustomobject]@{ Row = 25; B = $null; C = $null; D = "1.882"; E = "  -1.18%  " }
    [pscustomobject]@{ Row = 26; B = $null; C = $null; D = "151.15"; E = "  -0.82%  " }
    [pscustomobject]@{ Row = 27; B = $null; C = $null; D = "18.27"; E = "  -1.86%  " }
    [pscustomobject]@{ Row = 28; B = $null; C = $null; D = "2.085"; E = "  -3.65%  " }
    [pscustomobject]@{ Row = 29; B = $null; C = $null; D = "5.109"; E = "  -3.15%  " }
    [pscustomobject]@{ Row = 30; B = $null; C = $null; D = "116.13"; E = "  -1.35%  " }
    [pscustomobject]@{ Row = 31; B = $null; C = $null; D = "0.08908"; E = "  -0.22%  " }
    [pscustomobject]@{ Row = 32; B = $null; C = $null; D = $null; E = "  +0.51%  " }
    [pscustomobject]@{ Row = 33; B = $null; C = $null; D = "0.7383"; E = "  -2.80%  " }
    [pscustomobject]@{ Row = 34; B = $null; C = $null; D = "1.148"; E = "  -1.94%  " }
    [pscustomobject]@{ Row = 35; B = $null; C = $null; D = "4.456"; E = "  -2.01%  " }
    [pscustomobject]@{ Row = 36; B = $null; C = $null; D = "1.010"; E = "  -0.06%  " }
    [pscustomobject]@{ Row = 37; B = $null; C = $null; D = "2.499"; E = "  +3.37%  " }
    [pscustomobject]@{ Row = 38; B = $null; C = $null; D = $null; E = "  -1.89%  " }
    [pscustomobject]@{ Row = 39; B = $null; C = $null; D = "0.05235"; E = "  -1.93%  " }
    [pscustomobject]@{ Row = 40; B = $null; C = $null; D = "0.01914"; E = "  -2.43%  " }
    [pscustomobject]@{ Row = 41; B = $null; C = $null; D = "2.930"; E = "  -2.64%  " }
    [pscustomobject]@{ Row = 42; B = $null; C = $null; D = "7.157"; E = "  -2.40%  " }
    [pscustomobject]@{ Row = 43; B = $null; C = $null; D = "0.5215"; E = "  -2.48%  " }
    [pscustomobject]@{ Row = 44; B = $null; C = $null; D = "0.1627"; E = "  -2.08%  " }
    [pscustomobject]@{ Row = 45; B = $null; C = $null; D = "8.287"; E = "  -2.85%  " }
    [pscustomobject]@{ Row = 46; B = $null; C = $null; D = "0.4844"; E = "  -2.20%  " }
    [pscustomobject]@{ Row = 47; B = "PaxDollar"; C = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"; D = "1.010"; E = "  -0.12%  " }
    [pscustomobject]@{ Row = 48; B = "EnergySwap"; C = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D = "10.15"; E = "  -3.53%  " }
    [pscustomobject]@{ Row = 49; B = $null; C = $null; D = "103.61"; E = "  -0.26%  " }
    [pscustomobject]@{ Row = 50; B = $null; C = $null; D = "1.635"; E = "  -2.16%  " }
    [pscustomobject]@{ Row = 51; B = $null; C = $null; D = "0.06267"; E = "  -1.03%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.B) { $ws.Range("B$($u.Row)").Value = $u.B }
    if ($null -ne $u.C) { $ws.Range("C$($u.Row)").Value = $u.C }
    if ($null -ne $u.D) {
        # Force text so values like "1.012" / "0.4630" are not reinterpreted as numbers
        $cell = $ws.Range("D$($u.Row)")
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($null -ne $u.E) { $ws.Range("E$($u.Row)").Value = $u.E }
}
